$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated K (Strike#) column values for rows 2-18 in column G
$newValues = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 1
    6  = 2
    7  = 1
    8  = 2
    9  = 0
    10 = 2
    11 = 1
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 1
    17 = 1
    18 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
